$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Grok-4-Fast tag for GPT-5-mini Scenario No. 1',
    'Grok-4-Fast tag for GPT-5-mini Scenario No. 2',
    'Grok-4-Fast tag for GPT-5-mini Scenario No. 3',
    'Grok-4-Fast tag for GPT-5-mini Scenario No. 4',
    'Grok-4-Fast tag for GPT-4.1-nano Scenario No. 1',
    'Grok-4-Fast tag for GPT-4.1-nano Scenario No. 2',
    'Grok-4-Fast tag for GPT-4.1-nano Scenario No. 3',
    'Grok-4-Fast tag for GPT-4.1-nano Scenario No. 4',
    'Grok-4-Fast tag for GPT-5-chat Scenario No. 1',
    'Grok-4-Fast tag for GPT-5-chat Scenario No. 2',
    'Grok-4-Fast tag for GPT-5-chat Scenario No. 3',
    'Grok-4-Fast tag for GPT-5-chat Scenario No. 4',
    'Grok-4-Fast tag for ChatGPT-5-mini Scenario No. 1',
    'Grok-4-Fast tag for ChatGPT-5-mini Scenario No. 2',
    'Grok-4-Fast tag for ChatGPT-5-mini Scenario No. 3',
    'Grok-4-Fast tag for ChatGPT-5-mini Scenario No. 4',
    'Grok-4-Fast tag for Qwen3-32b Scenario No. 1',
    'Grok-4-Fast tag for Qwen3-32b Scenario No. 2',
    'Grok-4-Fast tag for Qwen3-32b Scenario No. 3',
    'Grok-4-Fast tag for Qwen3-32b Scenario No. 4',
    'Grok-4-Fast tag for Qwen3-14b Scenario No. 1',
    'Grok-4-Fast tag for Qwen3-14b Scenario No. 2',
    'Grok-4-Fast tag for Qwen3-14b Scenario No. 3',
    'Grok-4-Fast tag for Qwen3-14b Scenario No. 4',
    'Grok-4-Fast tag for Qwen3-8b Scenario No. 1',
    'Grok-4-Fast tag for Qwen3-8b Scenario No. 2',
    'Grok-4-Fast tag for Qwen3-8b Scenario No. 3',
    'Grok-4-Fast tag for Qwen3-8b Scenario No. 4',
    'Grok-4-Fast tag for Qwen2.5-7b Scenario No. 1',
    'Grok-4-Fast tag for Qwen2.5-7b Scenario No. 2',
    'Grok-4-Fast tag for Qwen2.5-7b Scenario No. 3',
    'Grok-4-Fast tag for Qwen2.5-7b Scenario No. 4',
    'Grok-4-Fast tag for Qwen2.5-72b Scenario No. 1',
    'Grok-4-Fast tag for Qwen2.5-72b Scenario No. 2',
    'Grok-4-Fast tag for Qwen2.5-72b Scenario No. 3',
    'Grok-4-Fast tag for Qwen2.5-72b Scenario No. 4',
    'Grok-4-Fast tag for Gemini-2.5-Flash Scenario No. 1',
    'Grok-4-Fast tag for Gemini-2.5-Flash Scenario No. 2',
    'Grok-4-Fast tag for Gemini-2.5-Flash Scenario No. 3',
    'Grok-4-Fast tag for Gemini-2.5-Flash Scenario No. 4',
    'Grok-4-Fast tag for Gemini-2.5-Flash-Lite Scenario No. 1',
    'Grok-4-Fast tag for Gemini-2.5-Flash-Lite Scenario No. 2',
    'Grok-4-Fast tag for Gemini-2.5-Flash-Lite Scenario No. 3',
    'Grok-4-Fast tag for Gemini-2.5-Flash-Lite Scenario No. 4',
    'Grok-4-Fast tag for Gemini-2.0-Flash-001 Scenario No. 1',
    'Grok-4-Fast tag for Gemini-2.0-Flash-001 Scenario No. 2',
    'Grok-4-Fast tag for Gemini-2.0-Flash-001 Scenario No. 3',
    'Grok-4-Fast tag for Gemini-2.0-Flash-001 Scenario No. 4',
    'Grok-4-Fast tag for Gemini-2.0-Flash-001-Lite Scenario No. 1',
    'Grok-4-Fast tag for Gemini-2.0-Flash-001-Lite Scenario No. 2',
    'Grok-4-Fast tag for Gemini-2.0-Flash-001-Lite Scenario No. 3',
    'Grok-4-Fast tag for Gemini-2.0-Flash-001-Lite Scenario No. 4',
    'Grok-4-Fast tag for Grok-4-Fast Scenario No. 1',
    'Grok-4-Fast tag for Grok-4-Fast Scenario No. 2',
    'Grok-4-Fast tag for Grok-4-Fast Scenario No. 3',
    'Grok-4-Fast tag for Grok-4-Fast Scenario No. 4',
    'Grok-4-Fast tag for Gemma-3-27b-It Scenario No. 1',
    'Grok-4-Fast tag for Gemma-3-27b-It Scenario No. 2',
    'Grok-4-Fast tag for Gemma-3-27b-It Scenario No. 3',
    'Grok-4-Fast tag for Gemma-3-27b-It Scenario No. 4',
    'Grok-4-Fast tag for Gemma-3-4b-It Scenario No. 1',
    'Grok-4-Fast tag for Gemma-3-4b-It Scenario No. 2',
    'Grok-4-Fast tag for Gemma-3-4b-It Scenario No. 3',
    'Grok-4-Fast tag for Gemma-3-4b-It Scenario No. 4',
    'Grok-4-Fast tag for Gemma-3-12b-It Scenario No. 1',
    'Grok-4-Fast tag for Gemma-3-12b-It Scenario No. 2',
    'Grok-4-Fast tag for Gemma-3-12b-It Scenario No. 3',
    'Grok-4-Fast tag for Gemma-3-12b-It Scenario No. 4',
    'Grok-4-Fast tag for Gemma-3n-2B Scenario No. 1',
    'Grok-4-Fast tag for Gemma-3n-2B Scenario No. 2',
    'Grok-4-Fast tag for Gemma-3n-2B Scenario No. 3',
    'Grok-4-Fast tag for Gemma-3n-2B Scenario No. 4',
    'Grok-4-Fast tag for Gemma-3n-4B Scenario No. 1',
    'Grok-4-Fast tag for Gemma-3n-4B Scenario No. 2',
    'Grok-4-Fast tag for Gemma-3n-4B Scenario No. 3',
    'Grok-4-Fast tag for Gemma-3n-4B Scenario No. 4',
    'Grok-4-Fast tag for Gemma-2-9b-It Scenario No. 1',
    'Grok-4-Fast tag for Gemma-2-9b-It Scenario No. 2',
    'Grok-4-Fast tag for Gemma-2-9b-It Scenario No. 3',
    'Grok-4-Fast tag for Gemma-2-9b-It Scenario No. 4',
    'Grok-4-Fast tag for Mistral-Small-3.2-24b Scenario No. 1',
    'Grok-4-Fast tag for Mistral-Small-3.2-24b Scenario No. 2',
    'Grok-4-Fast tag for Mistral-Small-3.2-24b Scenario No. 3',
    'Grok-4-Fast tag for Mistral-Small-3.2-24b Scenario No. 4',
    'Grok-4-Fast tag for Mistral-Small-24b-2501 Scenario No. 1',
    'Grok-4-Fast tag for Mistral-Small-24b-2501 Scenario No. 2',
    'Grok-4-Fast tag for Mistral-Small-24b-2501 Scenario No. 3',
    'Grok-4-Fast tag for Mistral-Small-24b-2501 Scenario No. 4',
    'Grok-4-Fast tag for Mistral-Medium-3 Scenario No. 1',
    'Grok-4-Fast tag for Mistral-Medium-3 Scenario No. 2',
    'Grok-4-Fast tag for Mistral-Medium-3 Scenario No. 3',
    'Grok-4-Fast tag for Mistral-Medium-3 Scenario No. 4',
    'Grok-4-Fast tag for Mistral-Small-3.1-24b Scenario No. 1',
    'Grok-4-Fast tag for Mistral-Small-3.1-24b Scenario No. 2',
    'Grok-4-Fast tag for Mistral-Small-3.1-24b Scenario No. 3',
    'Grok-4-Fast tag for Mistral-Small-3.1-24b Scenario No. 4',
    'Grok-4-Fast tag for Mistral-large-2 Scenario No. 1',
    'Grok-4-Fast tag for Mistral-large-2 Scenario No. 2',
    'Grok-4-Fast tag for Mistral-large-2 Scenario No. 3',
    'Grok-4-Fast tag for Mistral-large-2 Scenario No. 4',
    'Grok-4-Fast tag for Llama-3.3-70b Scenario No. 1',
    'Grok-4-Fast tag for Llama-3.3-70b Scenario No. 2',
    'Grok-4-Fast tag for Llama-3.3-70b Scenario No. 3',
    'Grok-4-Fast tag for Llama-3.3-70b Scenario No. 4',
    'Grok-4-Fast tag for Llama-3.1-8b Scenario No. 1',
    'Grok-4-Fast tag for Llama-3.1-8b Scenario No. 2',
    'Grok-4-Fast tag for Llama-3.1-8b Scenario No. 3',
    'Grok-4-Fast tag for Llama-3.1-8b Scenario No. 4',
    'Grok-4-Fast tag for Llama-3.1-405b Scenario No. 1',
    'Grok-4-Fast tag for Llama-3.1-405b Scenario No. 2',
    'Grok-4-Fast tag for Llama-3.1-405b Scenario No. 3',
    'Grok-4-Fast tag for Llama-3.1-405b Scenario No. 4',
    'Grok-4-Fast tag for Llama-3.2-90b Scenario No. 1',
    'Grok-4-Fast tag for Llama-3.2-90b Scenario No. 2',
    'Grok-4-Fast tag for Llama-3.2-90b Scenario No. 3',
    'Grok-4-Fast tag for Llama-3.2-90b Scenario No. 4',
    'Grok-4-Fast tag for Llama-3.2-1b Scenario No. 1',
    'Grok-4-Fast tag for Llama-3.2-1b Scenario No. 2',
    'Grok-4-Fast tag for Llama-3.2-1b Scenario No. 3',
    'Grok-4-Fast tag for Llama-3.2-1b Scenario No. 4',
    'Grok-4-Fast tag for Llama-3.2-3b Scenario No. 1',
    'Grok-4-Fast tag for Llama-3.2-3b Scenario No. 2',
    'Grok-4-Fast tag for Llama-3.2-3b Scenario No. 3',
    'Grok-4-Fast tag for Llama-3.2-3b Scenario No. 4',
    'Grok-4-Fast tag for Llama-4-Scout Scenario No. 1',
    'Grok-4-Fast tag for Llama-4-Scout Scenario No. 2',
    'Grok-4-Fast tag for Llama-4-Scout Scenario No. 3',
    'Grok-4-Fast tag for Llama-4-Scout Scenario No. 4',
    'Grok-4-Fast tag for Llama-4-Maverick Scenario No. 1',
    'Grok-4-Fast tag for Llama-4-Maverick Scenario No. 2',
    'Grok-4-Fast tag for Llama-4-Maverick Scenario No. 3',
    'Grok-4-Fast tag for Llama-4-Maverick Scenario No. 4',
    'Grok-4-Fast tag for Llama-3-8b Scenario No. 1',
    'Grok-4-Fast tag for Llama-3-8b Scenario No. 2',
    'Grok-4-Fast tag for Llama-3-8b Scenario No. 3',
    'Grok-4-Fast tag for Llama-3-8b Scenario No. 4',
    'Grok-4-Fast tag for Llama-3-70b Scenario No. 1',
    'Grok-4-Fast tag for Llama-3-70b Scenario No. 2',
    'Grok-4-Fast tag for Llama-3-70b Scenario No. 3',
    'Grok-4-Fast tag for Llama-3-70b Scenario No. 4',
    'Grok-4-Fast tag for Llama-3.3-8b Scenario No. 1',
    'Grok-4-Fast tag for Llama-3.3-8b Scenario No. 2',
    'Grok-4-Fast tag for Llama-3.3-8b Scenario No. 3',
    'Grok-4-Fast tag for Llama-3.3-8b Scenario No. 4',
    'Grok-4-Fast tag for Command-A_(Alt) Scenario No. 1',
    'Grok-4-Fast tag for Command-A_(Alt) Scenario No. 2',
    'Grok-4-Fast tag for Command-A_(Alt) Scenario No. 3',
    'Grok-4-Fast tag for Command-A_(Alt) Scenario No. 4',
    'Grok-4-Fast tag for Command-R-Plus-08-2024 Scenario No. 1',
    'Grok-4-Fast tag for Command-R-Plus-08-2024 Scenario No. 2',
    'Grok-4-Fast tag for Command-R-Plus-08-2024 Scenario No. 3',
    'Grok-4-Fast tag for Command-R-Plus-08-2024 Scenario No. 4',
    'Grok-4-Fast tag for Command-R-08-2024 Scenario No. 1',
    'Grok-4-Fast tag for Command-R-08-2024 Scenario No. 2',
    'Grok-4-Fast tag for Command-R-08-2024 Scenario No. 3',
    'Grok-4-Fast tag for Command-R-08-2024 Scenario No. 4',
    'Grok-4-Fast tag for Command-R7b Scenario No. 1',
    'Grok-4-Fast tag for Command-R7b Scenario No. 2',
    'Grok-4-Fast tag for Command-R7b Scenario No. 3',
    'Grok-4-Fast tag for Command-R7b Scenario No. 4',
    'Grok-4-Fast tag for DeepSeek-Chat-V3-0324 Scenario No. 1',
    'Grok-4-Fast tag for DeepSeek-Chat-V3-0324 Scenario No. 2',
    'Grok-4-Fast tag for DeepSeek-Chat-V3-0324 Scenario No. 3',
    'Grok-4-Fast tag for DeepSeek-Chat-V3-0324 Scenario No. 4',
    'Grok-4-Fast tag for DeepSeek-Chat-V3.1 Scenario No. 1',
    'Grok-4-Fast tag for DeepSeek-Chat-V3.1 Scenario No. 2',
    'Grok-4-Fast tag for DeepSeek-Chat-V3.1 Scenario No. 3',
    'Grok-4-Fast tag for DeepSeek-Chat-V3.1 Scenario No. 4',
    'Grok-4-Fast tag for DeepSeek-V3 Scenario No. 1',
    'Grok-4-Fast tag for DeepSeek-V3 Scenario No. 2',
    'Grok-4-Fast tag for DeepSeek-V3 Scenario No. 3',
    'Grok-4-Fast tag for DeepSeek-V3 Scenario No. 4',
    'Grok-4-Fast tag for Mistral-Small-24b-2501 Scenario No. 1',
    'Grok-4-Fast tag for Mistral-Small-24b-2501 Scenario No. 2',
    'Grok-4-Fast tag for Mistral-Small-24b-2501 Scenario No. 3',
    'Grok-4-Fast tag for Mistral-Small-24b-2501 Scenario No. 4',
    'Grok-4-Fast tag for Claude-Sonnet-4 Scenario No. 1',
    'Grok-4-Fast tag for Claude-Sonnet-4 Scenario No. 2',
    'Grok-4-Fast tag for Claude-Sonnet-4 Scenario No. 3',
    'Grok-4-Fast tag for Claude-Sonnet-4 Scenario No. 4',
    'Grok-4-Fast tag for Claude-3.5-Sonnet Scenario No. 1',
    'Grok-4-Fast tag for Claude-3.5-Sonnet Scenario No. 2',
    'Grok-4-Fast tag for Claude-3.5-Sonnet Scenario No. 3',
    'Grok-4-Fast tag for Claude-3.5-Sonnet Scenario No. 4',
    'Grok-4-Fast tag for Claude-Opus-4 Scenario No. 1',
    'Grok-4-Fast tag for Claude-Opus-4 Scenario No. 2',
    'Grok-4-Fast tag for Claude-Opus-4 Scenario No. 3',
    'Grok-4-Fast tag for Claude-Opus-4 Scenario No. 4',
    'Grok-4-Fast tag for Claude-4.5-Sonnet Scenario No. 1',
    'Grok-4-Fast tag for Claude-4.5-Sonnet Scenario No. 2',
    'Grok-4-Fast tag for Claude-4.5-Sonnet Scenario No. 3',
    'Grok-4-Fast tag for Claude-4.5-Sonnet Scenario No. 4',
    'Grok-4-Fast tag for Claude-4.5-Haiku Scenario No. 1',
    'Grok-4-Fast tag for Claude-4.5-Haiku Scenario No. 2',
    'Grok-4-Fast tag for Claude-4.5-Haiku Scenario No. 3',
    'Grok-4-Fast tag for Claude-4.5-Haiku Scenario No. 4',
    'Grok-4-Fast tag for Claude-3.5-Haiku Scenario No. 1',
    'Grok-4-Fast tag for Claude-3.5-Haiku Scenario No. 2',
    'Grok-4-Fast tag for Claude-3.5-Haiku Scenario No. 3',
    'Grok-4-Fast tag for Claude-3.5-Haiku Scenario No. 4',
    'Grok-4-Fast tag for Claude-3-Haiku Scenario No. 1',
    'Grok-4-Fast tag for Claude-3-Haiku Scenario No. 2',
    'Grok-4-Fast tag for Claude-3-Haiku Scenario No. 3',
    'Grok-4-Fast tag for Claude-3-Haiku Scenario No. 4',
    'Grok-4-Fast tag for Grok-3 Scenario No. 1',
    'Grok-4-Fast tag for Grok-3 Scenario No. 2',
    'Grok-4-Fast tag for Grok-3 Scenario No. 3',
    'Grok-4-Fast tag for Grok-3 Scenario No. 4',
    'Grok-4-Fast tag for Grok-4-Fast Scenario No. 1',
    'Grok-4-Fast tag for Grok-4-Fast Scenario No. 2',
    'Grok-4-Fast tag for Grok-4-Fast Scenario No. 3',
    'Grok-4-Fast tag for Grok-4-Fast Scenario No. 4',
    'Grok-4-Fast tag for Phi-4 Scenario No. 1',
    'Grok-4-Fast tag for Phi-4 Scenario No. 2',
    'Grok-4-Fast tag for Phi-4 Scenario No. 3',
    'Grok-4-Fast tag for Phi-4 Scenario No. 4',
    'Grok-4-Fast tag for Phi-3-mini Scenario No. 1',
    'Grok-4-Fast tag for Phi-3-mini Scenario No. 2',
    'Grok-4-Fast tag for Phi-3-mini Scenario No. 3',
    'Grok-4-Fast tag for Phi-3-mini Scenario No. 4',
    'Grok-4-Fast tag for Phi-3.5-mini Scenario No. 1',
    'Grok-4-Fast tag for Phi-3.5-mini Scenario No. 2',
    'Grok-4-Fast tag for Phi-3.5-mini Scenario No. 3',
    'Grok-4-Fast tag for Phi-3.5-mini Scenario No. 4',
    'Grok-4-Fast tag for Phi-3-medium Scenario No. 1',
    'Grok-4-Fast tag for Phi-3-medium Scenario No. 2',
    'Grok-4-Fast tag for Phi-3-medium Scenario No. 3',
    'Grok-4-Fast tag for Phi-3-medium Scenario No. 4',
    'GPT-5-mini tag for GPT-5-mini Scenario No. 1',
    'GPT-5-mini tag for GPT-5-mini Scenario No. 2',
    'GPT-5-mini tag for GPT-5-mini Scenario No. 3',
    'GPT-5-mini tag for GPT-5-mini Scenario No. 4',
    'GPT-5-mini tag for GPT-4.1-nano Scenario No. 1',
    'GPT-5-mini tag for GPT-4.1-nano Scenario No. 2',
    'GPT-5-mini tag for GPT-4.1-nano Scenario No. 3',
    'GPT-5-mini tag for GPT-4.1-nano Scenario No. 4',
    'GPT-5-mini tag for GPT-5-chat Scenario No. 1',
    'GPT-5-mini tag for GPT-5-chat Scenario No. 2',
    'GPT-5-mini tag for GPT-5-chat Scenario No. 3',
    'GPT-5-mini tag for GPT-5-chat Scenario No. 4',
    'GPT-5-mini tag for ChatGPT-5-mini Scenario No. 1',
    'GPT-5-mini tag for ChatGPT-5-mini Scenario No. 2',
    'GPT-5-mini tag for ChatGPT-5-mini Scenario No. 3',
    'GPT-5-mini tag for ChatGPT-5-mini Scenario No. 4',
    'GPT-5-mini tag for Qwen3-32b Scenario No. 1',
    'GPT-5-mini tag for Qwen3-32b Scenario No. 2',
    'GPT-5-mini tag for Qwen3-32b Scenario No. 3',
    'GPT-5-mini tag for Qwen3-32b Scenario No. 4',
    'GPT-5-mini tag for Qwen3-14b Scenario No. 1',
    'GPT-5-mini tag for Qwen3-14b Scenario No. 2',
    'GPT-5-mini tag for Qwen3-14b Scenario No. 3',
    'GPT-5-mini tag for Qwen3-14b Scenario No. 4',
    'GPT-5-mini tag for Qwen3-8b Scenario No. 1',
    'GPT-5-mini tag for Qwen3-8b Scenario No. 2',
    'GPT-5-mini tag for Qwen3-8b Scenario No. 3',
    'GPT-5-mini tag for Qwen3-8b Scenario No. 4',
    'GPT-5-mini tag for Qwen2.5-7b Scenario No. 1',
    'GPT-5-mini tag for Qwen2.5-7b Scenario No. 2',
    'GPT-5-mini tag for Qwen2.5-7b Scenario No. 3',
    'GPT-5-mini tag for Qwen2.5-7b Scenario No. 4',
    'GPT-5-mini tag for Qwen2.5-72b Scenario No. 1',
    'GPT-5-mini tag for Qwen2.5-72b Scenario No. 2',
    'GPT-5-mini tag for Qwen2.5-72b Scenario No. 3',
    'GPT-5-mini tag for Qwen2.5-72b Scenario No. 4',
    'GPT-5-mini tag for Gemini-2.5-Flash Scenario No. 1',
    'GPT-5-mini tag for Gemini-2.5-Flash Scenario No. 2',
    'GPT-5-mini tag for Gemini-2.5-Flash Scenario No. 3',
    'GPT-5-mini tag for Gemini-2.5-Flash Scenario No. 4',
    'GPT-5-mini tag for Gemini-2.5-Flash-Lite Scenario No. 1',
    'GPT-5-mini tag for Gemini-2.5-Flash-Lite Scenario No. 2',
    'GPT-5-mini tag for Gemini-2.5-Flash-Lite Scenario No. 3',
    'GPT-5-mini tag for Gemini-2.5-Flash-Lite Scenario No. 4',
    'GPT-5-mini tag for Gemini-2.0-Flash-001 Scenario No. 1',
    'GPT-5-mini tag for Gemini-2.0-Flash-001 Scenario No. 2',
    'GPT-5-mini tag for Gemini-2.0-Flash-001 Scenario No. 3',
    'GPT-5-mini tag for Gemini-2.0-Flash-001 Scenario No. 4',
    'GPT-5-mini tag for Gemini-2.0-Flash-001-Lite Scenario No. 1',
    'GPT-5-mini tag for Gemini-2.0-Flash-001-Lite Scenario No. 2',
    'GPT-5-mini tag for Gemini-2.0-Flash-001-Lite Scenario No. 3',
    'GPT-5-mini tag for Gemini-2.0-Flash-001-Lite Scenario No. 4',
    'GPT-5-mini tag for Grok-4-Fast Scenario No. 1',
    'GPT-5-mini tag for Grok-4-Fast Scenario No. 2',
    'GPT-5-mini tag for Grok-4-Fast Scenario No. 3',
    'GPT-5-mini tag for Grok-4-Fast Scenario No. 4',
    'GPT-5-mini tag for Gemma-3-27b-It Scenario No. 1',
    'GPT-5-mini tag for Gemma-3-27b-It Scenario No. 2',
    'GPT-5-mini tag for Gemma-3-27b-It Scenario No. 3',
    'GPT-5-mini tag for Gemma-3-27b-It Scenario No. 4',
    'GPT-5-mini tag for Gemma-3-4b-It Scenario No. 1',
    'GPT-5-mini tag for Gemma-3-4b-It Scenario No. 2',
    'GPT-5-mini tag for Gemma-3-4b-It Scenario No. 3',
    'GPT-5-mini tag for Gemma-3-4b-It Scenario No. 4',
    'GPT-5-mini tag for Gemma-3-12b-It Scenario No. 1',
    'GPT-5-mini tag for Gemma-3-12b-It Scenario No. 2',
    'GPT-5-mini tag for Gemma-3-12b-It Scenario No. 3',
    'GPT-5-mini tag for Gemma-3-12b-It Scenario No. 4',
    'GPT-5-mini tag for Gemma-3n-2B Scenario No. 1',
    'GPT-5-mini tag for Gemma-3n-2B Scenario No. 2',
    'GPT-5-mini tag for Gemma-3n-2B Scenario No. 3',
    'GPT-5-mini tag for Gemma-3n-2B Scenario No. 4',
    'GPT-5-mini tag for Gemma-3n-4B Scenario No. 1',
    'GPT-5-mini tag for Gemma-3n-4B Scenario No. 2',
    'GPT-5-mini tag for Gemma-3n-4B Scenario No. 3',
    'GPT-5-mini tag for Gemma-3n-4B Scenario No. 4',
    'GPT-5-mini tag for Gemma-2-9b-It Scenario No. 1',
    'GPT-5-mini tag for Gemma-2-9b-It Scenario No. 2',
    'GPT-5-mini tag for Gemma-2-9b-It Scenario No. 3',
    'GPT-5-mini tag for Gemma-2-9b-It Scenario No. 4',
    'GPT-5-mini tag for Mistral-Small-3.2-24b Scenario No. 1',
    'GPT-5-mini tag for Mistral-Small-3.2-24b Scenario No. 2',
    'GPT-5-mini tag for Mistral-Small-3.2-24b Scenario No. 3',
    'GPT-5-mini tag for Mistral-Small-3.2-24b Scenario No. 4',
    'GPT-5-mini tag for Mistral-Small-24b-2501 Scenario No. 1',
    'GPT-5-mini tag for Mistral-Small-24b-2501 Scenario No. 2',
    'GPT-5-mini tag for Mistral-Small-24b-2501 Scenario No. 3',
    'GPT-5-mini tag for Mistral-Small-24b-2501 Scenario No. 4',
    'GPT-5-mini tag for Mistral-Medium-3 Scenario No. 1',
    'GPT-5-mini tag for Mistral-Medium-3 Scenario No. 2',
    'GPT-5-mini tag for Mistral-Medium-3 Scenario No. 3',
    'GPT-5-mini tag for Mistral-Medium-3 Scenario No. 4',
    'GPT-5-mini tag for Mistral-Small-3.1-24b Scenario No. 1',
    'GPT-5-mini tag for Mistral-Small-3.1-24b Scenario No. 2',
    'GPT-5-mini tag for Mistral-Small-3.1-24b Scenario No. 3',
    'GPT-5-mini tag for Mistral-Small-3.1-24b Scenario No. 4',
    'GPT-5-mini tag for Mistral-large-2 Scenario No. 1',
    'GPT-5-mini tag for Mistral-large-2 Scenario No. 2',
    'GPT-5-mini tag for Mistral-large-2 Scenario No. 3',
    'GPT-5-mini tag for Mistral-large-2 Scenario No. 4',
    'GPT-5-mini tag for Llama-3.3-70b Scenario No. 1',
    'GPT-5-mini tag for Llama-3.3-70b Scenario No. 2',
    'GPT-5-mini tag for Llama-3.3-70b Scenario No. 3',
    'GPT-5-mini tag for Llama-3.3-70b Scenario No. 4',
    'GPT-5-mini tag for Llama-3.1-8b Scenario No. 1',
    'GPT-5-mini tag for Llama-3.1-8b Scenario No. 2',
    'GPT-5-mini tag for Llama-3.1-8b Scenario No. 3',
    'GPT-5-mini tag for Llama-3.1-8b Scenario No. 4',
    'GPT-5-mini tag for Llama-3.1-405b Scenario No. 1',
    'GPT-5-mini tag for Llama-3.1-405b Scenario No. 2',
    'GPT-5-mini tag for Llama-3.1-405b Scenario No. 3',
    'GPT-5-mini tag for Llama-3.1-405b Scenario No. 4',
    'GPT-5-mini tag for Llama-3.2-90b Scenario No. 1',
    'GPT-5-mini tag for Llama-3.2-90b Scenario No. 2',
    'GPT-5-mini tag for Llama-3.2-90b Scenario No. 3',
    'GPT-5-mini tag for Llama-3.2-90b Scenario No. 4',
    'GPT-5-mini tag for Llama-3.2-1b Scenario No. 1',
    'GPT-5-mini tag for Llama-3.2-1b Scenario No. 2',
    'GPT-5-mini tag for Llama-3.2-1b Scenario No. 3',
    'GPT-5-mini tag for Llama-3.2-1b Scenario No. 4',
    'GPT-5-mini tag for Llama-3.2-3b Scenario No. 1',
    'GPT-5-mini tag for Llama-3.2-3b Scenario No. 2',
    'GPT-5-mini tag for Llama-3.2-3b Scenario No. 3',
    'GPT-5-mini tag for Llama-3.2-3b Scenario No. 4',
    'GPT-5-mini tag for Llama-4-Scout Scenario No. 1',
    'GPT-5-mini tag for Llama-4-Scout Scenario No. 2',
    'GPT-5-mini tag for Llama-4-Scout Scenario No. 3',
    'GPT-5-mini tag for Llama-4-Scout Scenario No. 4',
    'GPT-5-mini tag for Llama-4-Maverick Scenario No. 1',
    'GPT-5-mini tag for Llama-4-Maverick Scenario No. 2',
    'GPT-5-mini tag for Llama-4-Maverick Scenario No. 3',
    'GPT-5-mini tag for Llama-4-Maverick Scenario No. 4',
    'GPT-5-mini tag for Llama-3-8b Scenario No. 1',
    'GPT-5-mini tag for Llama-3-8b Scenario No. 2',
    'GPT-5-mini tag for Llama-3-8b Scenario No. 3',
    'GPT-5-mini tag for Llama-3-8b Scenario No. 4',
    'GPT-5-mini tag for Llama-3-70b Scenario No. 1',
    'GPT-5-mini tag for Llama-3-70b Scenario No. 2',
    'GPT-5-mini tag for Llama-3-70b Scenario No. 3',
    'GPT-5-mini tag for Llama-3-70b Scenario No. 4',
    'GPT-5-mini tag for Llama-3.3-8b Scenario No. 1',
    'GPT-5-mini tag for Llama-3.3-8b Scenario No. 2',
    'GPT-5-mini tag for Llama-3.3-8b Scenario No. 3',
    'GPT-5-mini tag for Llama-3.3-8b Scenario No. 4',
    'GPT-5-mini tag for Command-A_(Alt) Scenario No. 1',
    'GPT-5-mini tag for Command-A_(Alt) Scenario No. 2',
    'GPT-5-mini tag for Command-A_(Alt) Scenario No. 3',
    'GPT-5-mini tag for Command-A_(Alt) Scenario No. 4',
    'GPT-5-mini tag for Command-R-Plus-08-2024 Scenario No. 1',
    'GPT-5-mini tag for Command-R-Plus-08-2024 Scenario No. 2',
    'GPT-5-mini tag for Command-R-Plus-08-2024 Scenario No. 3',
    'GPT-5-mini tag for Command-R-Plus-08-2024 Scenario No. 4',
    'GPT-5-mini tag for Command-R-08-2024 Scenario No. 1',
    'GPT-5-mini tag for Command-R-08-2024 Scenario No. 2',
    'GPT-5-mini tag for Command-R-08-2024 Scenario No. 3',
    'GPT-5-mini tag for Command-R-08-2024 Scenario No. 4',
    'GPT-5-mini tag for Command-R7b Scenario No. 1',
    'GPT-5-mini tag for Command-R7b Scenario No. 2',
    'GPT-5-mini tag for Command-R7b Scenario No. 3',
    'GPT-5-mini tag for Command-R7b Scenario No. 4',
    'GPT-5-mini tag for DeepSeek-Chat-V3-0324 Scenario No. 1',
    'GPT-5-mini tag for DeepSeek-Chat-V3-0324 Scenario No. 2',
    'GPT-5-mini tag for DeepSeek-Chat-V3-0324 Scenario No. 3',
    'GPT-5-mini tag for DeepSeek-Chat-V3-0324 Scenario No. 4',
    'GPT-5-mini tag for DeepSeek-Chat-V3.1 Scenario No. 1',
    'GPT-5-mini tag for DeepSeek-Chat-V3.1 Scenario No. 2',
    'GPT-5-mini tag for DeepSeek-Chat-V3.1 Scenario No. 3',
    'GPT-5-mini tag for DeepSeek-Chat-V3.1 Scenario No. 4',
    'GPT-5-mini tag for DeepSeek-V3 Scenario No. 1',
    'GPT-5-mini tag for DeepSeek-V3 Scenario No. 2',
    'GPT-5-mini tag for DeepSeek-V3 Scenario No. 3',
    'GPT-5-mini tag for DeepSeek-V3 Scenario No. 4',
    'GPT-5-mini tag for Mistral-Small-24b-2501 Scenario No. 1',
    'GPT-5-mini tag for Mistral-Small-24b-2501 Scenario No. 2',
    'GPT-5-mini tag for Mistral-Small-24b-2501 Scenario No. 3',
    'GPT-5-mini tag for Mistral-Small-24b-2501 Scenario No. 4',
    'GPT-5-mini tag for Claude-Sonnet-4 Scenario No. 1',
    'GPT-5-mini tag for Claude-Sonnet-4 Scenario No. 2',
    'GPT-5-mini tag for Claude-Sonnet-4 Scenario No. 3',
    'GPT-5-mini tag for Claude-Sonnet-4 Scenario No. 4',
    'GPT-5-mini tag for Claude-3.5-Sonnet Scenario No. 1',
    'GPT-5-mini tag for Claude-3.5-Sonnet Scenario No. 2',
    'GPT-5-mini tag for Claude-3.5-Sonnet Scenario No. 3',
    'GPT-5-mini tag for Claude-3.5-Sonnet Scenario No. 4',
    'GPT-5-mini tag for Claude-Opus-4 Scenario No. 1',
    'GPT-5-mini tag for Claude-Opus-4 Scenario No. 2',
    'GPT-5-mini tag for Claude-Opus-4 Scenario No. 3',
    'GPT-5-mini tag for Claude-Opus-4 Scenario No. 4',
    'GPT-5-mini tag for Claude-4.5-Sonnet Scenario No. 1',
    'GPT-5-mini tag for Claude-4.5-Sonnet Scenario No. 2',
    'GPT-5-mini tag for Claude-4.5-Sonnet Scenario No. 3',
    'GPT-5-mini tag for Claude-4.5-Sonnet Scenario No. 4',
    'GPT-5-mini tag for Claude-4.5-Haiku Scenario No. 1',
    'GPT-5-mini tag for Claude-4.5-Haiku Scenario No. 2',
    'GPT-5-mini tag for Claude-4.5-Haiku Scenario No. 3',
    'GPT-5-mini tag for Claude-4.5-Haiku Scenario No. 4',
    'GPT-5-mini tag for Claude-3.5-Haiku Scenario No. 1',
    'GPT-5-mini tag for Claude-3.5-Haiku Scenario No. 2',
    'GPT-5-mini tag for Claude-3.5-Haiku Scenario No. 3',
    'GPT-5-mini tag for Claude-3.5-Haiku Scenario No. 4',
    'GPT-5-mini tag for Claude-3-Haiku Scenario No. 1',
    'GPT-5-mini tag for Claude-3-Haiku Scenario No. 2',
    'GPT-5-mini tag for Claude-3-Haiku Scenario No. 3',
    'GPT-5-mini tag for Claude-3-Haiku Scenario No. 4',
    'GPT-5-mini tag for Grok-3 Scenario No. 1',
    'GPT-5-mini tag for Grok-3 Scenario No. 2',
    'GPT-5-mini tag for Grok-3 Scenario No. 3',
    'GPT-5-mini tag for Grok-3 Scenario No. 4',
    'GPT-5-mini tag for Grok-4-Fast Scenario No. 1',
    'GPT-5-mini tag for Grok-4-Fast Scenario No. 2',
    'GPT-5-mini tag for Grok-4-Fast Scenario No. 3',
    'GPT-5-mini tag for Grok-4-Fast Scenario No. 4',
    'GPT-5-mini tag for Phi-4 Scenario No. 1',
    'GPT-5-mini tag for Phi-4 Scenario No. 2',
    'GPT-5-mini tag for Phi-4 Scenario No. 3',
    'GPT-5-mini tag for Phi-4 Scenario No. 4',
    'GPT-5-mini tag for Phi-3-mini Scenario No. 1',
    'GPT-5-mini tag for Phi-3-mini Scenario No. 2',
    'GPT-5-mini tag for Phi-3-mini Scenario No. 3',
    'GPT-5-mini tag for Phi-3-mini Scenario No. 4',
    'GPT-5-mini tag for Phi-3.5-mini Scenario No. 1',
    'GPT-5-mini tag for Phi-3.5-mini Scenario No. 2',
    'GPT-5-mini tag for Phi-3.5-mini Scenario No. 3',
    'GPT-5-mini tag for Phi-3.5-mini Scenario No. 4',
    'GPT-5-mini tag for Phi-3-medium Scenario No. 1',
    'GPT-5-mini tag for Phi-3-medium Scenario No. 2',
    'GPT-5-mini tag for Phi-3-medium Scenario No. 3',
    'GPT-5-mini tag for Phi-3-medium Scenario No. 4',
    'Mistral-Small-24b-2501 tag for GPT-5-mini Scenario No. 1',
    'Mistral-Small-24b-2501 tag for GPT-5-mini Scenario No. 2',
    'Mistral-Small-24b-2501 tag for GPT-5-mini Scenario No. 3',
    'Mistral-Small-24b-2501 tag for GPT-5-mini Scenario No. 4',
    'Mistral-Small-24b-2501 tag for GPT-4.1-nano Scenario No. 1',
    'Mistral-Small-24b-2501 tag for GPT-4.1-nano Scenario No. 2',
    'Mistral-Small-24b-2501 tag for GPT-4.1-nano Scenario No. 3',
    'Mistral-Small-24b-2501 tag for GPT-4.1-nano Scenario No. 4',
    'Mistral-Small-24b-2501 tag for GPT-5-chat Scenario No. 1',
    'Mistral-Small-24b-2501 tag for GPT-5-chat Scenario No. 2',
    'Mistral-Small-24b-2501 tag for GPT-5-chat Scenario No. 3',
    'Mistral-Small-24b-2501 tag for GPT-5-chat Scenario No. 4',
    'Mistral-Small-24b-2501 tag for ChatGPT-5-mini Scenario No. 1',
    'Mistral-Small-24b-2501 tag for ChatGPT-5-mini Scenario No. 2',
    'Mistral-Small-24b-2501 tag for ChatGPT-5-mini Scenario No. 3',
    'Mistral-Small-24b-2501 tag for ChatGPT-5-mini Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Qwen3-32b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Qwen3-32b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Qwen3-32b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Qwen3-32b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Qwen3-14b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Qwen3-14b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Qwen3-14b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Qwen3-14b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Qwen3-8b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Qwen3-8b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Qwen3-8b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Qwen3-8b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Qwen2.5-7b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Qwen2.5-7b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Qwen2.5-7b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Qwen2.5-7b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Qwen2.5-72b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Qwen2.5-72b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Qwen2.5-72b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Qwen2.5-72b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemini-2.5-Flash Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemini-2.5-Flash Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemini-2.5-Flash Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemini-2.5-Flash Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemini-2.5-Flash-Lite Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemini-2.5-Flash-Lite Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemini-2.5-Flash-Lite Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemini-2.5-Flash-Lite Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemini-2.0-Flash-001 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemini-2.0-Flash-001 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemini-2.0-Flash-001 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemini-2.0-Flash-001 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemini-2.0-Flash-001-Lite Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemini-2.0-Flash-001-Lite Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemini-2.0-Flash-001-Lite Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemini-2.0-Flash-001-Lite Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Grok-4-Fast Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Grok-4-Fast Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Grok-4-Fast Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Grok-4-Fast Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemma-3-27b-It Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemma-3-27b-It Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemma-3-27b-It Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemma-3-27b-It Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemma-3-4b-It Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemma-3-4b-It Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemma-3-4b-It Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemma-3-4b-It Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemma-3-12b-It Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemma-3-12b-It Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemma-3-12b-It Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemma-3-12b-It Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemma-3n-2B Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemma-3n-2B Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemma-3n-2B Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemma-3n-2B Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemma-3n-4B Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemma-3n-4B Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemma-3n-4B Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemma-3n-4B Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemma-2-9b-It Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemma-2-9b-It Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemma-2-9b-It Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemma-2-9b-It Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Mistral-Small-3.2-24b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Mistral-Small-3.2-24b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Mistral-Small-3.2-24b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Mistral-Small-3.2-24b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Mistral-Small-24b-2501 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Mistral-Small-24b-2501 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Mistral-Small-24b-2501 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Mistral-Small-24b-2501 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Mistral-Medium-3 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Mistral-Medium-3 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Mistral-Medium-3 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Mistral-Medium-3 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Mistral-Small-3.1-24b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Mistral-Small-3.1-24b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Mistral-Small-3.1-24b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Mistral-Small-3.1-24b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Mistral-large-2 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Mistral-large-2 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Mistral-large-2 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Mistral-large-2 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-3.3-70b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-3.3-70b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-3.3-70b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-3.3-70b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-3.1-8b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-3.1-8b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-3.1-8b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-3.1-8b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-3.1-405b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-3.1-405b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-3.1-405b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-3.1-405b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-3.2-90b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-3.2-90b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-3.2-90b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-3.2-90b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-3.2-1b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-3.2-1b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-3.2-1b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-3.2-1b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-3.2-3b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-3.2-3b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-3.2-3b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-3.2-3b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-4-Scout Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-4-Scout Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-4-Scout Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-4-Scout Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-4-Maverick Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-4-Maverick Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-4-Maverick Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-4-Maverick Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-3-8b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-3-8b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-3-8b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-3-8b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-3-70b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-3-70b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-3-70b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-3-70b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Llama-3.3-8b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Llama-3.3-8b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Llama-3.3-8b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Llama-3.3-8b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Command-A_(Alt) Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Command-A_(Alt) Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Command-A_(Alt) Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Command-A_(Alt) Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Command-R-Plus-08-2024 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Command-R-Plus-08-2024 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Command-R-Plus-08-2024 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Command-R-Plus-08-2024 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Command-R-08-2024 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Command-R-08-2024 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Command-R-08-2024 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Command-R-08-2024 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Command-R7b Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Command-R7b Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Command-R7b Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Command-R7b Scenario No. 4',
    'Mistral-Small-24b-2501 tag for DeepSeek-Chat-V3-0324 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for DeepSeek-Chat-V3-0324 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for DeepSeek-Chat-V3-0324 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for DeepSeek-Chat-V3-0324 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for DeepSeek-Chat-V3.1 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for DeepSeek-Chat-V3.1 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for DeepSeek-Chat-V3.1 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for DeepSeek-Chat-V3.1 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for DeepSeek-V3 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for DeepSeek-V3 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for DeepSeek-V3 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for DeepSeek-V3 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Mistral-Small-24b-2501 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Mistral-Small-24b-2501 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Mistral-Small-24b-2501 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Mistral-Small-24b-2501 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Claude-Sonnet-4 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Claude-Sonnet-4 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Claude-Sonnet-4 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Claude-Sonnet-4 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Claude-3.5-Sonnet Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Claude-3.5-Sonnet Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Claude-3.5-Sonnet Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Claude-3.5-Sonnet Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Claude-Opus-4 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Claude-Opus-4 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Claude-Opus-4 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Claude-Opus-4 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Claude-4.5-Sonnet Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Claude-4.5-Sonnet Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Claude-4.5-Sonnet Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Claude-4.5-Sonnet Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Claude-4.5-Haiku Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Claude-4.5-Haiku Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Claude-4.5-Haiku Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Claude-4.5-Haiku Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Claude-3.5-Haiku Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Claude-3.5-Haiku Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Claude-3.5-Haiku Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Claude-3.5-Haiku Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Claude-3-Haiku Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Claude-3-Haiku Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Claude-3-Haiku Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Claude-3-Haiku Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Grok-3 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Grok-3 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Grok-3 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Grok-3 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Grok-4-Fast Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Grok-4-Fast Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Grok-4-Fast Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Grok-4-Fast Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Phi-4 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Phi-4 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Phi-4 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Phi-4 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Phi-3-mini Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Phi-3-mini Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Phi-3-mini Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Phi-3-mini Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Phi-3.5-mini Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Phi-3.5-mini Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Phi-3.5-mini Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Phi-3.5-mini Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Phi-3-medium Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Phi-3-medium Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Phi-3-medium Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Phi-3-medium Scenario No. 4',
    'Grok-4-Fast tag for GPT-4.1-mini Scenario No. 1',
    'Grok-4-Fast tag for GPT-4.1-mini Scenario No. 2',
    'Grok-4-Fast tag for GPT-4.1-mini Scenario No. 3',
    'Grok-4-Fast tag for GPT-4.1-mini Scenario No. 4',
    'Grok-4-Fast tag for GPT-4.1 Scenario No. 1',
    'Grok-4-Fast tag for GPT-4.1 Scenario No. 2',
    'Grok-4-Fast tag for GPT-4.1 Scenario No. 3',
    'Grok-4-Fast tag for GPT-4.1 Scenario No. 4',
    'Grok-4-Fast tag for GPT-5 Scenario No. 1',
    'Grok-4-Fast tag for GPT-5 Scenario No. 2',
    'Grok-4-Fast tag for GPT-5 Scenario No. 3',
    'Grok-4-Fast tag for GPT-5 Scenario No. 4',
    'Grok-4-Fast tag for GPT-5-mini Scenario No. 1',
    'Grok-4-Fast tag for GPT-5-mini Scenario No. 2',
    'Grok-4-Fast tag for GPT-5-mini Scenario No. 3',
    'Grok-4-Fast tag for GPT-5-mini Scenario No. 4',
    'Grok-4-Fast tag for GPT-5-nano Scenario No. 1',
    'Grok-4-Fast tag for GPT-5-nano Scenario No. 2',
    'Grok-4-Fast tag for GPT-5-nano Scenario No. 3',
    'Grok-4-Fast tag for GPT-5-nano Scenario No. 4',
    'GPT-5-mini tag for GPT-4.1-mini Scenario No. 1',
    'GPT-5-mini tag for GPT-4.1-mini Scenario No. 2',
    'GPT-5-mini tag for GPT-4.1-mini Scenario No. 3',
    'GPT-5-mini tag for GPT-4.1-mini Scenario No. 4',
    'GPT-5-mini tag for GPT-4.1 Scenario No. 1',
    'GPT-5-mini tag for GPT-4.1 Scenario No. 2',
    'GPT-5-mini tag for GPT-4.1 Scenario No. 3',
    'GPT-5-mini tag for GPT-4.1 Scenario No. 4',
    'GPT-5-mini tag for GPT-5 Scenario No. 1',
    'GPT-5-mini tag for GPT-5 Scenario No. 2',
    'GPT-5-mini tag for GPT-5 Scenario No. 3',
    'GPT-5-mini tag for GPT-5 Scenario No. 4',
    'GPT-5-mini tag for GPT-5-mini Scenario No. 1',
    'GPT-5-mini tag for GPT-5-mini Scenario No. 2',
    'GPT-5-mini tag for GPT-5-mini Scenario No. 3',
    'GPT-5-mini tag for GPT-5-mini Scenario No. 4',
    'GPT-5-mini tag for GPT-5-nano Scenario No. 1',
    'GPT-5-mini tag for GPT-5-nano Scenario No. 2',
    'GPT-5-mini tag for GPT-5-nano Scenario No. 3',
    'GPT-5-mini tag for GPT-5-nano Scenario No. 4',
    'Mistral-Small-24b-2501 tag for GPT-4.1-mini Scenario No. 1',
    'Mistral-Small-24b-2501 tag for GPT-4.1-mini Scenario No. 2',
    'Mistral-Small-24b-2501 tag for GPT-4.1-mini Scenario No. 3',
    'Mistral-Small-24b-2501 tag for GPT-4.1-mini Scenario No. 4',
    'Mistral-Small-24b-2501 tag for GPT-4.1 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for GPT-4.1 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for GPT-4.1 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for GPT-4.1 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for GPT-5 Scenario No. 1',
    'Mistral-Small-24b-2501 tag for GPT-5 Scenario No. 2',
    'Mistral-Small-24b-2501 tag for GPT-5 Scenario No. 3',
    'Mistral-Small-24b-2501 tag for GPT-5 Scenario No. 4',
    'Mistral-Small-24b-2501 tag for GPT-5-mini Scenario No. 1',
    'Mistral-Small-24b-2501 tag for GPT-5-mini Scenario No. 2',
    'Mistral-Small-24b-2501 tag for GPT-5-mini Scenario No. 3',
    'Mistral-Small-24b-2501 tag for GPT-5-mini Scenario No. 4',
    'Mistral-Small-24b-2501 tag for GPT-5-nano Scenario No. 1',
    'Mistral-Small-24b-2501 tag for GPT-5-nano Scenario No. 2',
    'Mistral-Small-24b-2501 tag for GPT-5-nano Scenario No. 3',
    'Mistral-Small-24b-2501 tag for GPT-5-nano Scenario No. 4',
    'Grok-4-Fast tag for Gemma-2-27b-It Scenario No. 1',
    'Grok-4-Fast tag for Gemma-2-27b-It Scenario No. 2',
    'Grok-4-Fast tag for Gemma-2-27b-It Scenario No. 3',
    'Grok-4-Fast tag for Gemma-2-27b-It Scenario No. 4',
    'GPT-5-mini tag for Gemma-2-27b-It Scenario No. 1',
    'GPT-5-mini tag for Gemma-2-27b-It Scenario No. 2',
    'GPT-5-mini tag for Gemma-2-27b-It Scenario No. 3',
    'GPT-5-mini tag for Gemma-2-27b-It Scenario No. 4',
    'Mistral-Small-24b-2501 tag for Gemma-2-27b-It Scenario No. 1',
    'Mistral-Small-24b-2501 tag for Gemma-2-27b-It Scenario No. 2',
    'Mistral-Small-24b-2501 tag for Gemma-2-27b-It Scenario No. 3',
    'Mistral-Small-24b-2501 tag for Gemma-2-27b-It Scenario No. 4'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $values[$i]
}
